$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the "datetimeFigureOut" date placeholder text on the slide
#    master and every slide layout (12/20/2024 -> 3/11/2025).
# ---------------------------------------------------------------------------
$newDate = "3/11/2025"

function Update-DatePlaceholder($shapes) {
    foreach ($shp in $shapes) {
        if ($shp.HasTextFrame) {
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            } catch {
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Rewrite the "3-day conference plus lodging" bullet so the pricing
#    callout reads "for less than $1500" instead of "for ~$1000".
# ---------------------------------------------------------------------------
$nbsp = [char]0xA0
$s = $p.Slides.Item(1)
foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "TextBox 7") {
        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($i = 1; $i -le $paraCount; $i++) {
            $para = $tr.Paragraphs($i, 1)
            if ($para.Text -like "*3-day conference plus lodging*") {
                $para.Text = "$nbsp $nbsp $nbsp $nbsp $nbsp- 3-day conference plus lodging "
                $para.InsertAfter("for less than `$1500$nbsp") | Out-Null
            }
        }
    }
}
